$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.027.84"
$ws.Range("E2").Value = "  -2.55%  "
$ws.Range("D3").Value = "3.519.83"
$ws.Range("E3").Value = "  -3.19%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.87"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.43%  "
$ws.Range("E7").Value = "  -0.44%  "
$ws.Range("D8").Value = "3.516.40"
$ws.Range("E8").Value = "  -3.09%  "
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("E10").Value = "  -4.52%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.75"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.16%  "
$ws.Range("E12").Value = "  -5.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "47.35"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.51%  "
$ws.Range("E14").Value = "  -3.34%  "
$ws.Range("D15").Value = "4.085.86"
$ws.Range("E15").Value = "  -3.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.43"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.88%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "612.06"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -9.49%  "
$ws.Range("D18").Value = "69.083.37"
$ws.Range("E18").Value = "  -2.53%  "
$ws.Range("D19").Value = "3.514.78"
$ws.Range("E19").Value = "  -3.32%  "
$ws.Range("E20").Value = "  -1.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.40"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.08"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.60%  "
$ws.Range("E23").Value = "  -6.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.75"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -8.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "96.46"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.75%  "
$ws.Range("E26").Value = "  -2.39%  "
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("E28").Value = "  -6.69%  "
$ws.Range("E29").Value = "  -6.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.57"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.51"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.69%  "
$ws.Range("E32").Value = "  -5.94%  "
$ws.Range("E33").Value = "  -5.36%  "
$ws.Range("E34").Value = "  -8.90%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "611.36"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.73"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.29%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.46"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -13.14%  "
$ws.Range("E38").Value = "  -4.87%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "57.14"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.44%  "
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0443"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.135"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.54%  "
$ws.Range("D43").Value = "3.389.65"
$ws.Range("E43").Value = "  -5.02%  "
$ws.Range("E44").Value = "  -5.84%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "32.75"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.89%  "
$ws.Range("E46").Value = "  -5.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.52"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.35%  "
$ws.Range("E48").Value = "  -6.34%  "
$ws.Range("E49").Value = "  -3.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "133.86"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.57"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +11.15%  "
